$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": append a new date column CT (19-sep) with 24 hourly
# price values, mirroring the existing CS column's header style.
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the header cell's formatting (bold, centered, bordered) from CS1
# onto CT1, then overwrite the value.
$wsPrix.Range("CS1").Copy()
$wsPrix.Range("CT1").PasteSpecial(-4122)
$wsPrix.Range("CT1").Value = "19-sep"

$prixValues = @{
    2  = 74.09
    3  = 65
    4  = 53.8
    5  = 46.93
    6  = 44.34
    7  = 42.13
    8  = 35.66
    9  = 81.19
    10 = 92.77
    11 = 49.14
    12 = 35.74
    13 = 13
    14 = 0.99
    15 = 0
    16 = 0
    17 = 11.97
    18 = 33.04
    19 = 65.87
    20 = 76
    21 = 108.14
    22 = 105.01
    23 = 85.59999999999999
    24 = 86.20999999999999
    25 = 79.14
}

foreach ($row in $prixValues.Keys) {
    $wsPrix.Range("CT$row").Value = $prixValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append row 95 for 2025-09-17.
# Use a formula-then-convert-to-values trick for the date cell so it is
# stored as plain text (matching the existing rows) instead of being
# auto-recognised as a date value/format.
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A95").Formula = '="2025-09-17"'
$wsGaz.Range("A95").Copy()
$wsGaz.Range("A95").PasteSpecial(-4163)
$wsGaz.Range("B95").Value = 31.725

# ---------------------------------------------------------------------
# Sheet "CO2": append row 95 for 2025-09-17.
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A95").Formula = '="2025-09-17"'
$wsCo2.Range("A95").Copy()
$wsCo2.Range("A95").PasteSpecial(-4163)
$wsCo2.Range("B95").Value = 76.5
